# Auto-generated edit script applying scheduled market-data refresh values
# to the Chocobo Profits workbook, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 25004700
$ws.Range("I74").Value = 50002000
$ws.Range("J74").Value = 7400
$ws.Range("K74").Value = 50002000
$ws.Range("L74").Value = 7400
$ws.Range("M74").Value = -50001064
$ws.Range("N74").Value = -9272
# Row 77
$ws.Range("H77").Value = 25004700
$ws.Range("I77").Value = 50002000
$ws.Range("J77").Value = 7400
$ws.Range("K77").Value = 250010000
$ws.Range("L77").Value = 37000
$ws.Range("M77").Value = -250005320
$ws.Range("N77").Value = -46360
# Row 100
$ws.Range("H100").Value = 14287303
$ws.Range("I100").Value = 15386250
$ws.Range("K100").Value = 15386250
$ws.Range("M100").Value = -15385709
# Row 103
$ws.Range("H103").Value = 740
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
# Row 113
$ws.Range("H113").Value = 6439.4
$ws.Range("I113").Value = 3748
$ws.Range("J113").Value = 7112.25
$ws.Range("K113").Value = 3748
$ws.Range("L113").Value = 7112.25
$ws.Range("M113").Value = -494
$ws.Range("N113").Value = -13620.25
# Row 125
$ws.Range("H125").Value = 1376.1428
$ws.Range("I125").Value = 1233.3334
$ws.Range("J125").Value = 1483.25
$ws.Range("K125").Value = 11100.0006
$ws.Range("L125").Value = 13349.25
$ws.Range("M125").Value = -8640.000599999999
$ws.Range("N125").Value = -18269.25
# Row 132
$ws.Range("H132").Value = 362204.3
$ws.Range("I132").Value = 5646.3184
$ws.Range("J132").Value = 1669583.6
$ws.Range("K132").Value = 16938.9552
$ws.Range("L132").Value = 5008750.800000001
$ws.Range("M132").Value = -14408.9552
$ws.Range("N132").Value = -5013810.800000001
# Row 138
$ws.Range("H138").Value = 4006.38
$ws.Range("I138").Value = 687.44446
$ws.Range("J138").Value = 5233.9316
$ws.Range("K138").Value = 2062.33338
$ws.Range("L138").Value = 15701.7948
$ws.Range("M138").Value = 3077.66662
$ws.Range("N138").Value = -25981.7948

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4517.9653
$ws.Range("I32").Value = 4512.4653
$ws.Range("K32").Value = 4512.4653
$ws.Range("M32").Value = -4225.4653
# Row 45
$ws.Range("H45").Value = 1098.4706
$ws.Range("I45").Value = 1107.5
$ws.Range("J45").Value = 1056.3334
$ws.Range("K45").Value = 1107.5
$ws.Range("L45").Value = 1056.3334
$ws.Range("M45").Value = -730.5
$ws.Range("N45").Value = -1810.3334
# Row 74
$ws.Range("H74").Value = 4864.16
$ws.Range("I74").Value = 5529.875
$ws.Range("J74").Value = 3680.6667
$ws.Range("K74").Value = 5529.875
$ws.Range("L74").Value = 3680.6667
$ws.Range("M74").Value = -4655.875
$ws.Range("N74").Value = -5428.6667
# Row 77
$ws.Range("H77").Value = 4864.16
$ws.Range("I77").Value = 5529.875
$ws.Range("J77").Value = 3680.6667
$ws.Range("K77").Value = 27649.375
$ws.Range("L77").Value = 18403.3335
$ws.Range("M77").Value = -23281.375
$ws.Range("N77").Value = -27139.3335
# Row 132
$ws.Range("H132").Value = 1938.975
$ws.Range("I132").Value = 941.96295
$ws.Range("J132").Value = 4009.6924
$ws.Range("K132").Value = 2825.88885
$ws.Range("L132").Value = 12029.0772
$ws.Range("M132").Value = -295.8888499999998
$ws.Range("N132").Value = -17089.0772

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1184.1818
$ws.Range("I107").Value = 1202.3
$ws.Range("J107").Value = 1003
$ws.Range("K107").Value = 1202.3
$ws.Range("L107").Value = 1003
$ws.Range("M107").Value = 717.7
$ws.Range("N107").Value = -4843
# Row 134
$ws.Range("H134").Value = 2294.3784
$ws.Range("I134").Value = 1476.963
$ws.Range("K134").Value = 4430.889
$ws.Range("M134").Value = -1895.889

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 16500
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 23000
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 23000
$ws.Range("M2").Value = -9887
$ws.Range("N2").Value = -23226
# Row 11
$ws.Range("H11").Value = 23670
$ws.Range("I11").Value = 350
$ws.Range("J11").Value = 29500
$ws.Range("K11").Value = 350
$ws.Range("L11").Value = 29500
$ws.Range("M11").Value = -210
$ws.Range("N11").Value = -29780
# Row 62
$ws.Range("H62").Value = 9000
$ws.Range("J62").Value = 9000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10248
# Row 65
$ws.Range("H65").Value = 9000
$ws.Range("J65").Value = 9000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51240
# Row 107
$ws.Range("H107").Value = 531.75
$ws.Range("I107").Value = 445.6
$ws.Range("J107").Value = 839.4286
$ws.Range("K107").Value = 445.6
$ws.Range("L107").Value = 839.4286
$ws.Range("M107").Value = 1474.4
$ws.Range("N107").Value = -4679.4286
# Row 132
$ws.Range("H132").Value = 2233.0645
$ws.Range("I132").Value = 1433.9166
$ws.Range("J132").Value = 4973
$ws.Range("K132").Value = 4301.7498
$ws.Range("L132").Value = 14919
$ws.Range("M132").Value = -1771.7498
$ws.Range("N132").Value = -19979
# Row 140
$ws.Range("H140").Value = 76634
$ws.Range("J140").Value = 76634
$ws.Range("L140").Value = 76634
$ws.Range("N140").Value = -86994

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 1696.5
$ws.Range("J13").Value = 2999.6667
$ws.Range("L13").Value = 8999.000100000001
$ws.Range("N13").Value = -9335.000100000001
# Row 70
$ws.Range("H70").Value = 2548.2942
$ws.Range("I70").Value = 1340.6666
$ws.Range("J70").Value = 2807.0715
$ws.Range("K70").Value = 4021.9998
$ws.Range("L70").Value = 8421.2145
$ws.Range("M70").Value = -3706.9998
$ws.Range("N70").Value = -9051.2145
# Row 73
$ws.Range("H73").Value = 2548.2942
$ws.Range("I73").Value = 1340.6666
$ws.Range("J73").Value = 2807.0715
$ws.Range("K73").Value = 4021.9998
$ws.Range("L73").Value = 8421.2145
$ws.Range("M73").Value = -2929.9998
$ws.Range("N73").Value = -10605.2145
# Row 81
$ws.Range("H81").Value = 1753.25
$ws.Range("I81").Value = 1006.5
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 3019.5
$ws.Range("L81").Value = 7500
$ws.Range("M81").Value = -1896.5
$ws.Range("N81").Value = -9746
# Row 82
$ws.Range("H82").Value = 5085.7144
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 5766.6665
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 17299.9995
$ws.Range("M82").Value = -2594
$ws.Range("N82").Value = -18111.9995
# Row 84
$ws.Range("H84").Value = 1753.25
$ws.Range("I84").Value = 1006.5
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 9058.5
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -3442.5
$ws.Range("N84").Value = -33732
# Row 85
$ws.Range("H85").Value = 5085.7144
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 5766.6665
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 17299.9995
$ws.Range("M85").Value = -1596
$ws.Range("N85").Value = -20107.9995
# Row 123
$ws.Range("H123").Value = 10000
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
# Row 131
$ws.Range("H131").Value = 8621514
$ws.Range("J131").Value = 867.94446
$ws.Range("L131").Value = 2603.83338
$ws.Range("N131").Value = -12683.83338
# Row 132
$ws.Range("H132").Value = 2567
$ws.Range("J132").Value = 2869.6155
$ws.Range("L132").Value = 25826.5395
$ws.Range("N132").Value = -30886.5395

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5943.859
$ws.Range("J70").Value = 7467.5713
$ws.Range("L70").Value = 7467.5713
$ws.Range("N70").Value = -8007.5713
# Row 73
$ws.Range("H73").Value = 5943.859
$ws.Range("J73").Value = 7467.5713
$ws.Range("L73").Value = 7467.5713
$ws.Range("N73").Value = -9339.5713
# Row 94
$ws.Range("H94").Value = 44599
$ws.Range("J94").Value = 44599
$ws.Range("L94").Value = 44599
$ws.Range("N94").Value = -45951
# Row 102
$ws.Range("H102").Value = 1620.125
$ws.Range("I102").Value = 1223.3572
$ws.Range("J102").Value = 4397.5
$ws.Range("K102").Value = 1223.3572
$ws.Range("L102").Value = 4397.5
$ws.Range("M102").Value = 398.6428000000001
$ws.Range("N102").Value = -7641.5
# Row 122
$ws.Range("H122").Value = 3379.0715
$ws.Range("I122").Value = 2485.1538
$ws.Range("K122").Value = 7455.4614
$ws.Range("M122").Value = -5005.4614
# Row 132
$ws.Range("H132").Value = 2135.525
$ws.Range("I132").Value = 1570.762
$ws.Range("K132").Value = 4712.286
$ws.Range("M132").Value = -2182.286

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3013.5
$ws.Range("I122").Value = 1667.25
$ws.Range("J122").Value = 7052.25
$ws.Range("K122").Value = 5001.75
$ws.Range("L122").Value = 21156.75
$ws.Range("M122").Value = -2551.75
$ws.Range("N122").Value = -26056.75
# Row 132
$ws.Range("H132").Value = 5564.227
$ws.Range("I132").Value = 1912
$ws.Range("K132").Value = 5736
$ws.Range("M132").Value = -3206

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 266.3846
$ws.Range("I107").Value = 271.66666
$ws.Range("J107").Value = 203
$ws.Range("K107").Value = 814.9999799999999
$ws.Range("L107").Value = 609
$ws.Range("M107").Value = 1105.00002
$ws.Range("N107").Value = -4449
# Row 122
$ws.Range("H122").Value = 3866.9375
$ws.Range("I122").Value = 2067
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 6201
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -3751
$ws.Range("N122").Value = -39900.001
# Row 132
$ws.Range("H132").Value = 15153608
$ws.Range("I132").Value = 1266.1428
$ws.Range("J132").Value = 41670204
$ws.Range("K132").Value = 3798.4284
$ws.Range("L132").Value = 125010612
$ws.Range("M132").Value = -1268.4284
$ws.Range("N132").Value = -125015672

